$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.715.15"
$ws.Range("E2").Value = "  -5.09%  "
$ws.Range("D3").Value = "3.017.79"
$ws.Range("E3").Value = "  -6.43%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.25"
$ws.Range("E5").Value = "  -2.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.30"
$ws.Range("E6").Value = "  -7.72%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "3.012.38"
$ws.Range("E8").Value = "  -6.54%  "
$ws.Range("E9").Value = "  -3.12%  "
$ws.Range("E10").Value = "  -7.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.16"
$ws.Range("E11").Value = "  -3.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.442"
$ws.Range("E12").Value = "  -3.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000225"
$ws.Range("E13").Value = "  -7.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.77"
$ws.Range("E14").Value = "  -8.46%  "
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").Value = "3.512.21"
$ws.Range("E16").Value = "  -6.59%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.011.65"
$ws.Range("E17").Value = "  -6.75%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "60.663.22"
$ws.Range("E18").Value = "  -5.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.42"
$ws.Range("E19").Value = "  -2.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "433.92"
$ws.Range("E20").Value = "  -7.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.15"
$ws.Range("E21").Value = "  -6.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.668"
$ws.Range("E22").Value = "  -5.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.04"
$ws.Range("E23").Value = "  -9.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.83"
$ws.Range("E24").Value = "  -5.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.64"
$ws.Range("E25").Value = "  -5.06%  "
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.58"
$ws.Range("E28").Value = "  -4.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.36"
$ws.Range("E29").Value = "  -7.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.93"
$ws.Range("E30").Value = "  -8.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.21"
$ws.Range("E31").Value = "  -10.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.43"
$ws.Range("E32").Value = "  -8.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0938"
$ws.Range("E33").Value = "  -9.59%  "
$ws.Range("E34").Value = "  -12.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.960"
$ws.Range("E35").Value = "  -8.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.62"
$ws.Range("E36").Value = "  -5.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "50.05"
$ws.Range("E37").Value = "  -3.44%  "
$ws.Range("D38").Value = "0.0₃0669"
$ws.Range("E38").Value = "  -9.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.49"
$ws.Range("E39").Value = "  +3.82%  "
$ws.Range("E40").Value = "  -9.29%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "388.80"
$ws.Range("E41").Value = "  -4.57%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.110"
$ws.Range("E42").Value = "  -3.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.51"
$ws.Range("E43").Value = "  -10.43%  "
$ws.Range("D44").Value = "2.666.93"
$ws.Range("E44").Value = "  -6.54%  "
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.237"
$ws.Range("E46").Value = "  -8.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.03"
$ws.Range("E47").Value = "  -6.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "119.23"
$ws.Range("E48").Value = "  -7.54%  "
$ws.Range("E49").Value = "  -4.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.81"
$ws.Range("E50").Value = "  -8.27%  "
$ws.Range("E51").Value = "  +1.84%  "
